$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Append a manual line break + a new bullet line to the
#    "4. Feature Specifications" Table-of-Contents entry (the
#    paragraph styled "List Number" near the top of the document).
#    NOTE: the same text also appears later as a "Heading 1" - that
#    occurrence must stay untouched.
# ------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq "4. Feature Specifications" -and $para.Style.NameLocal -eq "List Number") {
        $targetPara = $para
        break
    }
}

if ($targetPara -ne $null) {
    # Insert a line break right before the paragraph mark (End - 1).
    $endPos = $targetPara.Range.End
    $breakPoint = $d.Range($endPos - 1, $endPos - 1)
    $breakPoint.InsertAfter([char]11)

    # Insert the new bullet text right after the break (also before
    # the paragraph mark, which has shifted one character to the right).
    $endPos2 = $targetPara.Range.End
    $textPoint = $d.Range($endPos2 - 1, $endPos2 - 1)
    $bulletText = [char]0x2022 + " Factorial Operation: Calculates the factorial of a number, with input validation and error handling for negative numbers or non-integer values."
    $textPoint.InsertAfter($bulletText)
}

# ------------------------------------------------------------------
# 2. Add a new "Factorial" / "factorial" row at the end of the
#    Option/Function reference table (identified by its header row
#    text so the script does not depend on table ordering).
# ------------------------------------------------------------------
$optionTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $h1 = $tbl.Rows.Item(1).Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7)
    $h2 = $tbl.Rows.Item(1).Cells.Item(2).Range.Text.TrimEnd([char]13, [char]7)
    if ($h1 -eq "Option" -and $h2 -eq "Function") {
        $optionTable = $tbl
        break
    }
}

if ($optionTable -ne $null) {
    $newRow = $optionTable.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = "Factorial"
    $newRow.Cells.Item(2).Range.Text = "factorial"
}
